$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date column (A2:A3) -------------------------------------------------
# "2025-12-05" -> "2025-12-06" for both data rows. The source value is
# stored as text (not a real date), so force a text number format before
# typing it in (otherwise Excel auto-converts it to a date serial), then
# restore the default "Normal" style so the cell ends up identical to how
# it started, just with new text.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-06"
$ws.Range("A2").Style = "Normal"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-12-06"
$ws.Range("A3").Style = "Normal"

# --- Row 2 becomes the Archer Aviation (ACHR) record ---------------------
$ws.Range("B2").Value = "Archer Aviation Inc."
$ws.Range("C2").Value = "ACHR"
$ws.Range("D2").Value = 8.82
$ws.Range("E2").Value = 60.8
$ws.Range("F2").Value = 13.29
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 46
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 57.5
$ws.Range("N2").Value = 51.54219175917372

# --- Row 3 becomes the Joby Aviation (JOBY) record ------------------------
$ws.Range("B3").Value = "Joby Aviation, Inc."
$ws.Range("C3").Value = "JOBY"
$ws.Range("D3").Value = 15.4
$ws.Range("E3").Value = 57.7
$ws.Range("F3").Value = 6.74
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 57.5
$ws.Range("N3").Value = 51.54219175917372
